$d = $word.ActiveDocument

$d.Content.Find.Execute("74-5=", $true, $true, $false, $false, $false, $true, 1, $false, "15+3=", 2) | Out-Null
$d.Content.Find.Execute("99-23=", $true, $true, $false, $false, $false, $true, 1, $false, "34+37=", 2) | Out-Null
$d.Content.Find.Execute("39-16=", $true, $true, $false, $false, $false, $true, 1, $false, "91-22=", 2) | Out-Null
$d.Content.Find.Execute("98-7=", $true, $true, $false, $false, $false, $true, 1, $false, "96+2=", 2) | Out-Null
$d.Content.Find.Execute("45+20=", $true, $true, $false, $false, $false, $true, 1, $false, "80-52=", 2) | Out-Null
$d.Content.Find.Execute("18+71=", $true, $true, $false, $false, $false, $true, 1, $false, "74-74=", 2) | Out-Null
$d.Content.Find.Execute("74-60=", $true, $true, $false, $false, $false, $true, 1, $false, "8+80=", 2) | Out-Null
$d.Content.Find.Execute("76-3=", $true, $true, $false, $false, $false, $true, 1, $false, "9-6=", 2) | Out-Null
$d.Content.Find.Execute("46+47=", $true, $true, $false, $false, $false, $true, 1, $false, "70-48=", 2) | Out-Null
$d.Content.Find.Execute("70-16=", $true, $true, $false, $false, $false, $true, 1, $false, "5+31=", 2) | Out-Null
$d.Content.Find.Execute("47+25=", $true, $true, $false, $false, $false, $true, 1, $false, "30-7=", 2) | Out-Null
$d.Content.Find.Execute("97-16=", $true, $true, $false, $false, $false, $true, 1, $false, "58+17=", 2) | Out-Null
$d.Content.Find.Execute("43-30=", $true, $true, $false, $false, $false, $true, 1, $false, "33+53=", 2) | Out-Null
$d.Content.Find.Execute("92-83=", $true, $true, $false, $false, $false, $true, 1, $false, "26+12=", 2) | Out-Null
$d.Content.Find.Execute("25+56=", $true, $true, $false, $false, $false, $true, 1, $false, "16+67=", 2) | Out-Null
$d.Content.Find.Execute("32-10=", $true, $true, $false, $false, $false, $true, 1, $false, "96-58=", 2) | Out-Null
$d.Content.Find.Execute("36+12=", $true, $true, $false, $false, $false, $true, 1, $false, "34+42=", 2) | Out-Null
$d.Content.Find.Execute("40-17=", $true, $true, $false, $false, $false, $true, 1, $false, "53+17=", 2) | Out-Null
$d.Content.Find.Execute("1+87=", $true, $true, $false, $false, $false, $true, 1, $false, "42+37=", 2) | Out-Null
$d.Content.Find.Execute("72-69=", $true, $true, $false, $false, $false, $true, 1, $false, "24-0=", 2) | Out-Null
$d.Content.Find.Execute("19+4=", $true, $true, $false, $false, $false, $true, 1, $false, "6+92=", 2) | Out-Null
$d.Content.Find.Execute("15+66=", $true, $true, $false, $false, $false, $true, 1, $false, "83-24=", 2) | Out-Null
$d.Content.Find.Execute("76-27=", $true, $true, $false, $false, $false, $true, 1, $false, "56-36=", 2) | Out-Null
$d.Content.Find.Execute("79-26=", $true, $true, $false, $false, $false, $true, 1, $false, "65+14=", 2) | Out-Null
$d.Content.Find.Execute("68-0=", $true, $true, $false, $false, $false, $true, 1, $false, "6+33=", 2) | Out-Null
$d.Content.Find.Execute("70-22=", $true, $true, $false, $false, $false, $true, 1, $false, "48-31=", 2) | Out-Null
$d.Content.Find.Execute("61+12=", $true, $true, $false, $false, $false, $true, 1, $false, "16+38=", 2) | Out-Null
$d.Content.Find.Execute("85-84=", $true, $true, $false, $false, $false, $true, 1, $false, "33+24=", 2) | Out-Null
$d.Content.Find.Execute("88-42=", $true, $true, $false, $false, $false, $true, 1, $false, "31+21=", 2) | Out-Null
$d.Content.Find.Execute("10+73=", $true, $true, $false, $false, $false, $true, 1, $false, "80-52=", 2) | Out-Null
$d.Content.Find.Execute("31-15=", $true, $true, $false, $false, $false, $true, 1, $false, "94-52=", 2) | Out-Null
$d.Content.Find.Execute("20+74=", $true, $true, $false, $false, $false, $true, 1, $false, "59-41=", 2) | Out-Null
$d.Content.Find.Execute("72-19=", $true, $true, $false, $false, $false, $true, 1, $false, "95-56=", 2) | Out-Null
$d.Content.Find.Execute("19+72=", $true, $true, $false, $false, $false, $true, 1, $false, "10+32=", 2) | Out-Null
$d.Content.Find.Execute("91-53=", $true, $true, $false, $false, $false, $true, 1, $false, "76-63=", 2) | Out-Null
$d.Content.Find.Execute("40+7=", $true, $true, $false, $false, $false, $true, 1, $false, "55-38=", 2) | Out-Null
$d.Content.Find.Execute("85+3=", $true, $true, $false, $false, $false, $true, 1, $false, "39+30=", 2) | Out-Null
$d.Content.Find.Execute("22+20=", $true, $true, $false, $false, $false, $true, 1, $false, "20+63=", 2) | Out-Null
$d.Content.Find.Execute("89+8=", $true, $true, $false, $false, $false, $true, 1, $false, "52+42=", 2) | Out-Null
$d.Content.Find.Execute("83+16=", $true, $true, $false, $false, $false, $true, 1, $false, "38+26=", 2) | Out-Null
$d.Content.Find.Execute("10+76=", $true, $true, $false, $false, $false, $true, 1, $false, "89+9=", 2) | Out-Null
$d.Content.Find.Execute("46-5=", $true, $true, $false, $false, $false, $true, 1, $false, "65-51=", 2) | Out-Null
$d.Content.Find.Execute("96-89=", $true, $true, $false, $false, $false, $true, 1, $false, "94-81=", 2) | Out-Null
$d.Content.Find.Execute("70-3=", $true, $true, $false, $false, $false, $true, 1, $false, "29+46=", 2) | Out-Null
$d.Content.Find.Execute("86-20=", $true, $true, $false, $false, $false, $true, 1, $false, "40+33=", 2) | Out-Null
$d.Content.Find.Execute("37+42=", $true, $true, $false, $false, $false, $true, 1, $false, "64+24=", 2) | Out-Null
$d.Content.Find.Execute("36+40=", $true, $true, $false, $false, $false, $true, 1, $false, "94-51=", 2) | Out-Null
$d.Content.Find.Execute("41+41=", $true, $true, $false, $false, $false, $true, 1, $false, "69-12=", 2) | Out-Null
$d.Content.Find.Execute("36+42=", $true, $true, $false, $false, $false, $true, 1, $false, "69-60=", 2) | Out-Null
$d.Content.Find.Execute("73+9=", $true, $true, $false, $false, $false, $true, 1, $false, "95-5=", 2) | Out-Null
$d.Content.Find.Execute("81+13=", $true, $true, $false, $false, $false, $true, 1, $false, "10+77=", 2) | Out-Null
$d.Content.Find.Execute("10+70=", $true, $true, $false, $false, $false, $true, 1, $false, "85-68=", 2) | Out-Null
$d.Content.Find.Execute("29+0=", $true, $true, $false, $false, $false, $true, 1, $false, "10+3=", 2) | Out-Null
$d.Content.Find.Execute("31-27=", $true, $true, $false, $false, $false, $true, 1, $false, "3+59=", 2) | Out-Null
$d.Content.Find.Execute("72-35=", $true, $true, $false, $false, $false, $true, 1, $false, "61+9=", 2) | Out-Null
$d.Content.Find.Execute("45+7=", $true, $true, $false, $false, $false, $true, 1, $false, "8+48=", 2) | Out-Null
$d.Content.Find.Execute("90-5=", $true, $true, $false, $false, $false, $true, 1, $false, "52+26=", 2) | Out-Null
$d.Content.Find.Execute("83-49=", $true, $true, $false, $false, $false, $true, 1, $false, "69+18=", 2) | Out-Null
$d.Content.Find.Execute("71-39=", $true, $true, $false, $false, $false, $true, 1, $false, "26-20=", 2) | Out-Null
$d.Content.Find.Execute("2+69=", $true, $true, $false, $false, $false, $true, 1, $false, "69-48=", 2) | Out-Null
$d.Content.Find.Execute("51+31=", $true, $true, $false, $false, $false, $true, 1, $false, "56-32=", 2) | Out-Null
$d.Content.Find.Execute("43+33=", $true, $true, $false, $false, $false, $true, 1, $false, "8+10=", 2) | Out-Null
$d.Content.Find.Execute("59-36=", $true, $true, $false, $false, $false, $true, 1, $false, "12-10=", 2) | Out-Null
$d.Content.Find.Execute("96-18=", $true, $true, $false, $false, $false, $true, 1, $false, "89-62=", 2) | Out-Null
$d.Content.Find.Execute("82-73=", $true, $true, $false, $false, $false, $true, 1, $false, "90+2=", 2) | Out-Null
$d.Content.Find.Execute("8+12=", $true, $true, $false, $false, $false, $true, 1, $false, "72+21=", 2) | Out-Null
$d.Content.Find.Execute("68-27=", $true, $true, $false, $false, $false, $true, 1, $false, "70+8=", 2) | Out-Null
$d.Content.Find.Execute("53-50=", $true, $true, $false, $false, $false, $true, 1, $false, "93-86=", 2) | Out-Null
$d.Content.Find.Execute("31+54=", $true, $true, $false, $false, $false, $true, 1, $false, "94-36=", 2) | Out-Null
$d.Content.Find.Execute("20+70=", $true, $true, $false, $false, $false, $true, 1, $false, "39+16=", 2) | Out-Null
$d.Content.Find.Execute("20+1=", $true, $true, $false, $false, $false, $true, 1, $false, "98-64=", 2) | Out-Null
$d.Content.Find.Execute("8+91=", $true, $true, $false, $false, $false, $true, 1, $false, "28+41=", 2) | Out-Null
$d.Content.Find.Execute("86+0=", $true, $true, $false, $false, $false, $true, 1, $false, "27+36=", 2) | Out-Null
$d.Content.Find.Execute("64-32=", $true, $true, $false, $false, $false, $true, 1, $false, "85-19=", 2) | Out-Null
$d.Content.Find.Execute("82-68=", $true, $true, $false, $false, $false, $true, 1, $false, "15+25=", 2) | Out-Null
$d.Content.Find.Execute("16+34=", $true, $true, $false, $false, $false, $true, 1, $false, "65-47=", 2) | Out-Null
$d.Content.Find.Execute("96-50=", $true, $true, $false, $false, $false, $true, 1, $false, "96-84=", 2) | Out-Null
$d.Content.Find.Execute("84-76=", $true, $true, $false, $false, $false, $true, 1, $false, "76+16=", 2) | Out-Null
$d.Content.Find.Execute("58+35=", $true, $true, $false, $false, $false, $true, 1, $false, "25-4=", 2) | Out-Null
$d.Content.Find.Execute("90-60=", $true, $true, $false, $false, $false, $true, 1, $false, "55+18=", 2) | Out-Null
$d.Content.Find.Execute("51-19=", $true, $true, $false, $false, $false, $true, 1, $false, "15+32=", 2) | Out-Null
$d.Content.Find.Execute("6+11=", $true, $true, $false, $false, $false, $true, 1, $false, "72-58=", 2) | Out-Null
$d.Content.Find.Execute("78-40=", $true, $true, $false, $false, $false, $true, 1, $false, "39-23=", 2) | Out-Null
$d.Content.Find.Execute("69-21=", $true, $true, $false, $false, $false, $true, 1, $false, "58+36=", 2) | Out-Null
$d.Content.Find.Execute("27+1=", $true, $true, $false, $false, $false, $true, 1, $false, "97-91=", 2) | Out-Null
$d.Content.Find.Execute("38-16=", $true, $true, $false, $false, $false, $true, 1, $false, "22-16=", 2) | Out-Null
$d.Content.Find.Execute("82-2=", $true, $true, $false, $false, $false, $true, 1, $false, "65+3=", 2) | Out-Null
$d.Content.Find.Execute("99-93=", $true, $true, $false, $false, $false, $true, 1, $false, "96-39=", 2) | Out-Null
$d.Content.Find.Execute("31+11=", $true, $true, $false, $false, $false, $true, 1, $false, "14+82=", 2) | Out-Null
$d.Content.Find.Execute("94-61=", $true, $true, $false, $false, $false, $true, 1, $false, "2+70=", 2) | Out-Null
$d.Content.Find.Execute("22-17=", $true, $true, $false, $false, $false, $true, 1, $false, "0+91=", 2) | Out-Null
$d.Content.Find.Execute("31-30=", $true, $true, $false, $false, $false, $true, 1, $false, "5+19=", 2) | Out-Null
$d.Content.Find.Execute("26-21=", $true, $true, $false, $false, $false, $true, 1, $false, "26+16=", 2) | Out-Null
$d.Content.Find.Execute("44+41=", $true, $true, $false, $false, $false, $true, 1, $false, "24+73=", 2) | Out-Null
$d.Content.Find.Execute("18-5=", $true, $true, $false, $false, $false, $true, 1, $false, "27+70=", 2) | Out-Null
$d.Content.Find.Execute("40+20=", $true, $true, $false, $false, $false, $true, 1, $false, "38-34=", 2) | Out-Null
$d.Content.Find.Execute("87-53=", $true, $true, $false, $false, $false, $true, 1, $false, "49+25=", 2) | Out-Null
$d.Content.Find.Execute("36+50=", $true, $true, $false, $false, $false, $true, 1, $false, "67+2=", 2) | Out-Null
$d.Content.Find.Execute("22+11=", $true, $true, $false, $false, $false, $true, 1, $false, "3+38=", 2) | Out-Null
$d.Content.Find.Execute("37-17=", $true, $true, $false, $false, $false, $true, 1, $false, "67+30=", 2) | Out-Null
